# The document has one section whose header/footer each carry two logo
# images (a Pearson Edexcel logo in the footers, a BTec logo in the
# headers). Word exposes the "default" header/footer as Item(1) and the
# "first page" header/footer as Item(2). Per the authored change, the
# picture display-names are swapped between the two sibling image files
# that were embedded in the package (image1.png <-> image2.png for the
# Pearson logo, image1.jpg <-> image2.jpg for the BTec logo), while the
# pictures' visible content/position is left untouched.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footers: Pearson Edexcel logo -> rename image1.png to image2.png
$ftrDefault = $sec.Footers.Item(1)
if ($ftrDefault.Exists -and $ftrDefault.Range.InlineShapes.Count -gt 0) {
    $ftrDefault.Range.InlineShapes.Item(1).Name = "image2.png"
}

$ftrFirst = $sec.Footers.Item(2)
if ($ftrFirst.Exists -and $ftrFirst.Range.InlineShapes.Count -gt 0) {
    $ftrFirst.Range.InlineShapes.Item(1).Name = "image2.png"
}

# Headers: BTec logo -> rename image2.jpg to image1.jpg
$hdrDefault = $sec.Headers.Item(1)
if ($hdrDefault.Exists -and $hdrDefault.Range.InlineShapes.Count -gt 0) {
    $hdrDefault.Range.InlineShapes.Item(1).Name = "image1.jpg"
}

$hdrFirst = $sec.Headers.Item(2)
if ($hdrFirst.Exists -and $hdrFirst.Range.InlineShapes.Count -gt 0) {
    $hdrFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"
}
